# Update column G ("K") values on the active worksheet for rows 2-29.
# This mirrors a regeneration of save_data that now uses K (strike count)
# computed/written from s_vals instead of the old Strike# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2 through 29 (in order).
$newKValues = @(5, 2, 2, 3, 1, 5, 4, 5, 6, 4, 5, 3, 5, 3, 2, 4, 6, 3, 1, 3, 4, 6, 5, 3, 6, 4, 5, 1)

$startRow = 2
for ($i = 0; $i -lt $newKValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("G$row").Value = $newKValues[$i]
}
